# Updated cryptos list — refresh Price (col D) and Volume(1h) (col E) values
# for rows 2-51 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = "<new price text>"; E = "<new volume text>" }
# (rows with no 'D' key only had their Volume(1h) column change)
$updates = [ordered]@{
    2  = @{ D = "29.024.35";       E = "  -0.16%  " }
    3  = @{ D = "1.833.45";        E = "  -0.02%  " }
    4  = @{ D = "0.9984";          E = "  -0.12%  " }
    5  = @{ D = "242.22";          E = "  +0.12%  " }
    6  = @{ D = "0.6266";          E = "  -5.00%  " }
    7  = @{ D = "0.9994";          E = "  -0.24%  " }
    8  = @{ D = "0.07616";         E = "  +3.62%  " }
    9  = @{                        E = "  -0.36%  " }
    10 = @{                        E = "  -2.07%  " }
    11 = @{                        E = "  +0.14%  " }
    12 = @{ D = "1.835.15";        E = "  -0.18%  " }
    13 = @{ D = "4.954";           E = "  -0.65%  " }
    14 = @{ D = "0.6656";          E = "  -0.40%  " }
    15 = @{ D = "0.00001018";      E = "  +17.42%  " }
    16 = @{ D = "82.73";           E = "  +1.15%  " }
    17 = @{ D = "6.052";           E = "  -1.07%  " }
    18 = @{ D = "29.036.39";       E = "  -0.19%  " }
    19 = @{ D = "226.57";          E = "  +0.84%  " }
    20 = @{ D = "12.36";           E = "  -0.79%  " }
    21 = @{ D = "0.9987";          E = "  -0.24%  " }
    22 = @{ D = "7.183";           E = "  +0.78%  " }
    23 = @{                        E = "  -0.22%  " }
    24 = @{ D = "158.22";          E = "  +0.06%  " }
    25 = @{ D = "8.497";           E = "  -0.37%  " }
    26 = @{                        E = "  -0.73%  " }
    27 = @{ D = "17.90";           E = "  -0.12%  " }
    28 = @{ D = "1.490";           E = "  -1.29%  " }
    29 = @{ D = "4.109";           E = "  -0.17%  " }
    30 = @{ D = "4.017";           E = "  -0.15%  " }
    31 = @{                        E = "  -1.17%  " }
    32 = @{ D = "0.05217";         E = "  -3.18%  " }
    33 = @{ D = "1.844";           E = "  +0.26%  " }
    34 = @{ D = "0.7359";          E = "  -0.98%  " }
    35 = @{ D = "1.140";           E = "  -1.32%  " }
    36 = @{ D = "2.706";           E = "  +1.96%  " }
    37 = @{ D = "1.241.66";        E = "  -4.10%  " }
    38 = @{ D = "2.757";           E = "  -0.19%  " }
    39 = @{                        E = "  -0.35%  " }
    40 = @{ D = "6.345";           E = "  +0.09%  " }
    41 = @{ D = "0.8964";          E = "  -0.73%  " }
    42 = @{ D = "0.9995";          E = "  -0.12%  " }
    43 = @{ D = "101.58";          E = "  -1.69%  " }
    44 = @{ D = "1.982.57";        E = "  -0.33%  " }
    45 = @{ D = "0.00000000123";   E = "  -1.00%  " }
    46 = @{                        E = "  -0.51%  " }
    47 = @{ D = "0.5104";          E = "  -0.51%  " }
    48 = @{                        E = "  +0.92%  " }
    49 = @{ D = "8.866";           E = "  +1.80%  " }
    50 = @{ D = "1.644";           E = "  -5.74%  " }
    51 = @{ D = "0.05749";         E = "  -1.67%  " }
}

function Test-NumericLike([string]$s) {
    return $s -match '^[+-]?(\d+\.?\d*|\.\d+)$'
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    if ($vals.Contains('D')) {
        $dCell = $ws.Cells.Item($row, 4)
        $newD = $vals['D']
        # The Price column holds plain text (e.g. "29.024.35", "0.9984").
        # Some of these look like genuine numbers to Excel's auto-detection,
        # which would silently convert them to numeric values on assignment.
        # Force the cell to Text format first so the literal string sticks,
        # exactly like the values already sitting in the other Price cells.
        if (Test-NumericLike $newD) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $newD
    }

    if ($vals.Contains('E')) {
        $ws.Cells.Item($row, 5).Value = $vals['E']
    }
}
